# Updated cryptos list on Tue Aug 20 20:56:55 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "59.337.04"
$ws.Range("E2").Value = "  +0.60%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.584.98"
$ws.Range("E3").Value = "  -0.80%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "571.71"

# Row 6 - Solana
$ws.Range("D6").Value = "143.66"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.10%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +0.35%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.595.79"
$ws.Range("E9").Value = "  -0.95%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").Value = "  -1.67%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +3.16%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +10.85%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +2.63%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.045.96"
$ws.Range("E14").Value = "  -0.61%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "59.322.83"
$ws.Range("E15").Value = "  +0.64%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "22.56"
$ws.Range("E16").Value = "  +7.70%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +3.63%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.595.42"
$ws.Range("E18").Value = "  -0.60%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.50%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "336.08"
$ws.Range("E20").Value = "  -0.75%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "10.26"
$ws.Range("E21").Value = "  +1.24%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.72%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.02%  "

# Row 24 - Litecoin
$ws.Range("E24").Value = "  -3.34%  "

# Row 25 - Polygon
$ws.Range("D25").Value = "0.455"
$ws.Range("E25").Value = "  +6.38%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  +0.11%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -0.05%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  +1.42%  "

# Row 29 - PEPE
$ws.Range("D29").Value = "0.0₃0782"
$ws.Range("E29").Value = "  +3.30%  "

# Row 30 - USDe
$ws.Range("E30").Value = "  +0.03%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.30%  "

# Row 32 - Aptos
$ws.Range("E32").Value = "  +1.76%  "

# Row 33 - Monero
$ws.Range("D33").Value = "158.68"
$ws.Range("E33").Value = "  +3.17%  "

# Row 34 - EthereumClassic
$ws.Range("E34").Value = "  +0.37%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "4.06"
$ws.Range("E35").Value = "  +2.65%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +1.48%  "

# Row 37 <-> Row 38 swap (Fetch.AI / SuiNetwork) with updated values
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "0.876"
$ws.Range("E37").Value = "  -2.27%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "0.881"
$ws.Range("E38").Value = "  -0.68%  "

# Row 39 - OKB
$ws.Range("E39").Value = "  +0.50%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +2.68%  "

# Row 41 <-> Row 42 swap (Bittensor / Filecoin) with updated values
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.69"
$ws.Range("E41").Value = "  +2.03%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "294.83"
$ws.Range("E42").Value = "  +3.73%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.12%  "

# Row 44 - Stellar
$ws.Range("E44").Value = "  +2.81%  "

# Row 45 - Mantle
$ws.Range("D45").Value = "0.593"
$ws.Range("E45").Value = "  -1.49%  "

# Row 46 - EnergySwap
$ws.Range("D46").Value = "19.32"
$ws.Range("E46").Value = "  +2.45%  "

# Row 47 - Hedera
$ws.Range("E47").Value = "  +0.27%  "

# Row 48 - WhiteBITCoin
$ws.Range("E48").Value = "  -0.03%  "

# Row 49 - Aave
$ws.Range("D49").Value = "125.50"
$ws.Range("E49").Value = "  +6.76%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  +1.78%  "

# Row 51 - Maker
$ws.Range("D51").Value = "1.953.07"
$ws.Range("E51").Value = "  +0.19%  "
